# Realestate Update resale numbers 2023-06-17 18:18
# Appends a new data row (row 53) to the CityResaleNum sheet, mirroring
# the structure of the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 53

# --- Text columns (A-D) ---
# "Date" and "Week" look like a date / number to Excel's auto-detection,
# so force them to remain literal text (matching the existing rows) and
# then strip the "Text" number-format style that gets attached, so the
# new row's cells end up unstyled just like the rest of the sheet.

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-17"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "18:17:24"

$ws.Cells.Item($row, 3).Value = "Saturday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "24"
$ws.Cells.Item($row, 4).ClearFormats()

# --- Numeric columns (E-T) ---
$ws.Cells.Item($row, 5).Value = 122037
$ws.Cells.Item($row, 6).Value = 133282
$ws.Cells.Item($row, 7).Value = 162216
$ws.Cells.Item($row, 8).Value = 133199
$ws.Cells.Item($row, 9).Value = 177193
$ws.Cells.Item($row, 10).Value = 114767
$ws.Cells.Item($row, 11).Value = 201200
$ws.Cells.Item($row, 12).Value = 224946
$ws.Cells.Item($row, 13).Value = 175039
$ws.Cells.Item($row, 14).Value = 103365
$ws.Cells.Item($row, 15).Value = 39167
$ws.Cells.Item($row, 16).Value = 34019
$ws.Cells.Item($row, 17).Value = 51743
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36785
$ws.Cells.Item($row, 20).Value = -1
